{"js": "// Adjusted risk calc formula\n// The single table in the document has, for watershed \"Artlish\", one row\n// per limiting-factor (LF) description with Rank / Total Risk / Current\n// Risk / Future Risk columns. The edit re-sorts five adjacent rows\n// (0-indexed rows 8-12, i.e. the LF68/LF40/LF5/LF10/LF53 rows) by moving\n// the LF68 description (and its risk numbers) down past LF40/LF5/LF10/LF53,\n// which bumps the other four rows' Rank/Total Risk/Current Risk/Future Risk\n// values. We express this as per-cell text updates, addressed by\n// row/column index, so it is robust regardless of which literal strings\n// happen to repeat elsewhere in the table.\n\nconst table = context.document.body.tables.getFirst();\n\n// row index (0-based, header is row 0) -> column index (0-based) -> new text\nconst updates = [\n  // LF68 row becomes the LF40 row; Current/Future Risk swap M<->L\n  [8, 1, \"LF40: Mortality or fitness reduction due to frequent and higher peak flows causing flushing\"],\n  [8, 4, \"L\"],\n  [8, 5, \"M\"],\n\n  // LF40 row becomes the LF5 row (Rank/Total/Current/Future unchanged)\n  [9, 1, \"LF5: Mortality or fitness reduction due to competition with invasive species\"],\n\n  // LF5 row becomes the LF10 row; Rank/Total Risk/Future Risk shift\n  [10, 1, \"LF10: Mortality or fitness reduction of wild fish due to competition with hatchery fish or aquaculture escapees for spawning locations or mates\"],\n  [10, 2, \"9\"],\n  [10, 3, \"4\"],\n  [10, 5, \"L\"],\n\n  // LF10 row becomes the LF53 row; same Rank/Total Risk/Future Risk shift\n  [11, 1, \"LF53: Mortality or fitness reduction due to increased frequency and magnitude of algal blooms\"],\n  [11, 2, \"9\"],\n  [11, 3, \"4\"],\n  [11, 5, \"L\"],\n\n  // LF53 row becomes the LF68 row (moved to the end of this block) with\n  // updated Rank/Total Risk/Current Risk\n  [12, 1, \"LF68: Mortality or fitness reduction due to a reduction in natural (wild) genetic influence. This is measured by the stray rate (pHOSstray) into the system, or by the frequency and magnitude of direct transplanting.\"],\n  [12, 2, \"12\"],\n  [12, 3, \"3\"],\n  [12, 4, \"M\"],\n];\n\nfor (const [row, col, text] of updates) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Adjusted risk calc formula\n# The single table in the document has, for watershed \"Artlish\", one row\n# per limiting-factor (LF) description with Rank / Total Risk / Current\n# Risk / Future Risk columns. The edit re-sorts five adjacent rows\n# (Word 1-based rows 9-13, i.e. the LF68/LF40/LF5/LF10/LF53 rows) by moving\n# the LF68 description (and its risk numbers) down past LF40/LF5/LF10/LF53,\n# which bumps the other four rows' Rank/Total Risk/Current Risk/Future Risk\n# values. We express this as per-cell text updates, addressed by\n# row/column index (Word COM is 1-based for both Table.Cell(row, col) and\n# columns), so it is robust regardless of which literal strings happen to\n# repeat elsewhere in the table.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# Word COM rows/cols are 1-based; header is row 1, so data row (0-based n)\n# above is Word row (n + 1).\n\n# LF68 row (Word row 9) becomes the LF40 row; Current/Future Risk swap M<->L\n$tbl.Cell(9, 2).Range.Text = \"LF40: Mortality or fitness reduction due to frequent and higher peak flows causing flushing\"\n$tbl.Cell(9, 5).Range.Text = \"L\"\n$tbl.Cell(9, 6).Range.Text = \"M\"\n\n# LF40 row (Word row 10) becomes the LF5 row (Rank/Total/Current/Future unchanged)\n$tbl.Cell(10, 2).Range.Text = \"LF5: Mortality or fitness reduction due to competition with invasive species\"\n\n# LF5 row (Word row 11) becomes the LF10 row; Rank/Total Risk/Future Risk shift\n$tbl.Cell(11, 2).Range.Text = \"LF10: Mortality or fitness reduction of wild fish due to competition with hatchery fish or aquaculture escapees for spawning locations or mates\"\n$tbl.Cell(11, 3).Range.Text = \"9\"\n$tbl.Cell(11, 4).Range.Text = \"4\"\n$tbl.Cell(11, 6).Range.Text = \"L\"\n\n# LF10 row (Word row 12) becomes the LF53 row; same Rank/Total Risk/Future Risk shift\n$tbl.Cell(12, 2).Range.Text = \"LF53: Mortality or fitness reduction due to increased frequency and magnitude of algal blooms\"\n$tbl.Cell(12, 3).Range.Text = \"9\"\n$tbl.Cell(12, 4).Range.Text = \"4\"\n$tbl.Cell(12, 6).Range.Text = \"L\"\n\n# LF53 row (Word row 13) becomes the LF68 row (moved to the end of this\n# block) with updated Rank/Total Risk/Current Risk\n$tbl.Cell(13, 2).Range.Text = \"LF68: Mortality or fitness reduction due to a reduction in natural (wild) genetic influence. This is measured by the stray rate (pHOSstray) into the system, or by the frequency and magnitude of direct transplanting.\"\n$tbl.Cell(13, 3).Range.Text = \"12\"\n$tbl.Cell(13, 4).Range.Text = \"3\"\n$tbl.Cell(13, 5).Range.Text = \"M\"\n"}
